# 10-OCT-2024 One question solved on binary tree
# Adds a new "Binary Tree" section (rows 68-70) to Sheet1, mirroring the
# layout already used for the other topic sections in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 68: section header "Binary Tree" -------------------------------
# Reuse the formatting of an existing bold/size-20 header cell (B2) via
# copy/paste-special so no new font/style entries get minted.
$ws.Range("B2").Copy()
$ws.Cells.Item(68, 2).PasteSpecial(-4122)
$ws.Cells.Item(68, 2).Value2 = "                                                    Binary Tree"
$ws.Rows.Item(68).RowHeight = 31.5

# --- Row 69: problem entry ----------------------------------------------
# Reuse the date-format styling of an existing "s=3" cell (A67).
$ws.Range("A67").Copy()
$ws.Cells.Item(69, 1).PasteSpecial(-4122)
$ws.Cells.Item(69, 1).Value = "10/15/2024"
$ws.Cells.Item(69, 2).Value2 = "Subtree of another tree"
$ws.Cells.Item(69, 3).Value2 = "Recursion"
$ws.Cells.Item(69, 7).Value2 = "Easy"

# --- Row 70: notes row ----------------------------------------------------
# Write the text cells first so the new shared strings land in the same
# order as the source workbook; the literal "15/10/2024" text (kept in the
# same date-format style, but not parsed as a date) is added last.
$ws.Cells.Item(70, 2).Value2 = "calculate max depth of binary tree "
$ws.Cells.Item(70, 8).Value2 = "I have solved it in first attempt"
$ws.Range("A67").Copy()
$ws.Cells.Item(70, 1).PasteSpecial(-4122)
$ws.Cells.Item(70, 1).Value = "15/10/2024"

$excel.CutCopyMode = $false

# --- Update the view: scroll so row 53 is at the top and select B59 -----
$ws.Activate()
$ws.Range("B59").Select()
$excel.ActiveWindow.ScrollRow = 53
